$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario")

# Update stock (Existencia) for Paracetamol 250mg (row 2)
$ws.Range("D2").Value = 14

# Remove the "Talco para pies" row (old row 7) entirely, shifting rows 8-16 up.
$ws.Rows(7).Delete() | Out-Null

# The old row 8 (Tukol D Adulto 125ml) is now row 7 after the deletion above;
# update its stock (Existencia).
$ws.Range("D7").Value = 21

# The old row 10 (Fabe Naproxeno Paracetamol 15 tabletas) is now row 9 after the
# deletion above; replace it with the new "10 tabletas" product data.
$ws.Range("B9").Value = "Fabe Naproxeno Paracetamol 10 tabletas"
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = 25
